$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.760.76'
$ws.Range("D3").Value = '2.477.68'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.88'
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.22'
$ws.Range("E6").Value = '  +1.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E10").Value = '  +10.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '33.26'
$ws.Range("E11").Value = '  +3.10%  '
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").Value = '2.860.60'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.76'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '2.475.73'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.793'
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = '41.701.34'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.29'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.33'
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.83'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.76'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.74'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.83'
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.16'
$ws.Range("E30").Value = '  +2.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.44'
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("E35").Value = '  +1.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.51'
$ws.Range("E36").Value = '  +1.17%  '
$ws.Range("E37").Value = '  +6.22%  '
$ws.Range("E38").Value = '  +2.44%  '
$ws.Range("E39").Value = '  +1.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.103'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("E42").Value = '  +11.07%  '
$ws.Range("D43").Value = '1.993.23'
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.99'
$ws.Range("E45").Value = '  +3.10%  '
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.46'
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("D48").Value = '2.717.56'
$ws.Range("E48").Value = '  +0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.64'
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.24'
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.35'
$ws.Range("E51").Value = '  +0.37%  '
